$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "51.188.56"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.952.73"
$ws.Range("E3").Value = "  -0.64%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws.Range("D5") "382.86"
$ws.Range("E5").Value = "  +0.75%  "
Set-TextValue $ws.Range("D6") "102.65"
$ws.Range("E6").Value = "  -2.25%  "
Set-TextValue $ws.Range("D7") "0.539"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  +0.07%  "
Set-TextValue $ws.Range("D9") "0.590"
$ws.Range("E9").Value = "  -1.80%  "
Set-TextValue $ws.Range("D10") "36.63"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "3.427.50"
$ws.Range("E13").Value = "  +0.00%  "
Set-TextValue $ws.Range("D14") "18.08"
$ws.Range("E14").Value = "  -2.95%  "
Set-TextValue $ws.Range("D15") "7.41"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "2.961.88"
$ws.Range("E16").Value = "  +0.79%  "
Set-TextValue $ws.Range("D17") "0.988"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "51.178.24"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -7.84%  "
$ws.Range("E20").Value = "  -3.83%  "
Set-TextValue $ws.Range("D21") "12.59"
$ws.Range("E21").Value = "  -4.57%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  -0.49%  "
Set-TextValue $ws.Range("D23") "68.58"
$ws.Range("E23").Value = "  -0.22%  "
Set-TextValue $ws.Range("D24") "262.49"
$ws.Range("E24").Value = "  -0.45%  "
Set-TextValue $ws.Range("D25") "2.93"
$ws.Range("E25").Value = "  +4.71%  "
Set-TextValue $ws.Range("D26") "8.39"
$ws.Range("E26").Value = "  +12.89%  "
Set-TextValue $ws.Range("D27") "7.61"
$ws.Range("E27").Value = "  +1.97%  "
Set-TextValue $ws.Range("D28") "0.170"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +8.51%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -1.04%  "
Set-TextValue $ws.Range("D33") "34.29"
$ws.Range("E33").Value = "  -0.55%  "
Set-TextValue $ws.Range("D34") "0.0456"
$ws.Range("E34").Value = "  +3.92%  "
Set-TextValue $ws.Range("D35") "50.57"
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  +0.13%  "
Set-TextValue $ws.Range("D38") "2.99"
$ws.Range("E38").Value = "  -2.58%  "
Set-TextValue $ws.Range("D39") "16.87"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  -3.58%  "
Set-TextValue $ws.Range("D43") "121.17"
$ws.Range("E43").Value = "  -2.32%  "
Set-TextValue $ws.Range("D44") "21.44"
$ws.Range("E44").Value = "  -3.14%  "
$ws.Range("E45").Value = "  -1.25%  "
Set-TextValue $ws.Range("D46") "2.38"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("E47").Value = "  -3.95%  "
Set-TextValue $ws.Range("D48") "3.25"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "2.013.23"
$ws.Range("E49").Value = "  -1.04%  "
Set-TextValue $ws.Range("D50") "0.0349"
$ws.Range("E50").Value = "  +5.26%  "
$ws.Range("E51").Value = "  -1.30%  "
